$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-11 (A..G): Colaborador_id, Colaborador_nome, Departamento,
# Motivo_da_ausência, Horas_de_ausência, Data_da_ausência, Salário
$data = @(
    @(40357, "Esther Campos",        "Recursos Humanos",       "Viagem de negócios", 5, 45095, 8841.360000000001),
    @(91225, "Yasmin Ferreira",      "Jurídico",                "Problemas pessoais", 8, 45090, 2928.71),
    @(86177, "Nicolas Nascimento",   "Engenharia",               "Outros",             1, 45104, 7986.31),
    @(67099, "Sarah Monteiro",       "Atendimento ao Cliente",  "Consulta médica",    7, 45106, 9132.700000000001),
    @(97873, "João Pedro Moreira",   "P&D",                      "Problemas pessoais", 1, 45097, 4891.29),
    @(73934, "João Felipe Aragão",   "Jurídico",                "Outros",             4, 45088, 11649.57),
    @(91985, "João Gabriel da Rosa", "TI",                       "Outros",             4, 45082, 6319.63),
    @(18576, "Rebeca Costa",         "Marketing",                "Viagem de negócios", 3, 45084, 3267.52),
    @(94968, "Igor da Luz",          "Vendas",                   "Outros",             3, 45083, 5799.63),
    @(68838, "Srta. Bianca Dias",    "Jurídico",                "Consulta médica",    7, 45078, 6577.96)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $row++
}
